$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '$2b$10$zKZEiP/EXLr87BTtyuxcfugfzFoFxwHVNqrEvmL6hlGVrLHgcXGKe'
$ws.Range("C3").Value = '$2b$10$SxqbS/CYGZG/ar/CodD4JuJd7s/zg49VHOTeLvR9SARnp5bcwdIVu'
$ws.Range("C4").Value = '$2b$10$yjCF95I3vOITNzd4x.IEveyyyjkDWaXH72ANWB8l2e71YjdeQKYDm'
$ws.Range("C5").Value = '$2b$10$ECCIP0D81HKXpwe2yIiFL.zX1ZCIWVHD7yH5MUdzkTdgculSkPm42'
$ws.Range("C6").Value = '$2b$10$EZ3LrmJp1kcJO0LCI0eBE.3FWtVB6t2Dc1C3QWW5iJ4V2RUf5miES'
$ws.Range("C7").Value = '$2b$10$ggMSKZESkhbp6hOmQy7xzePh6i9FNMWPKx/UgWknCD4TnWWb6k/8S'
$ws.Range("C8").Value = '$2b$10$C47l/SUnIgKcsKaP4yWXLuFLueLuN5gI/QYnplP95RwPYlC7Z009.'
$ws.Range("C9").Value = '$2b$10$ScxEvRDrVzIazf4j/7IysOX6WUq1L85XvyZs7UNmb28qa3GSvWA5y'
$ws.Range("C10").Value = '$2b$10$XaQPC/JR8lE15Sn.BvxbcOFfDdVE2cfXtlNXMfPf5o9H.wr2zyzl.'
$ws.Range("C11").Value = '$2b$10$W1cs3CzREK5YPdR2.crjTe5BPnm.xbS86j1sqSBO2Uil.X432LBqW'
$ws.Range("C12").Value = '$2b$10$bNAAOJ/wxUrTyf42GeZYcuCbeDeqpXcgE6P2p7gqDFjqa74glO.d2'
$ws.Range("C13").Value = '$2b$10$fJKoPt9t/ObJj6O8mj/zUu1tQESU6Tw3.fAxAsNke5BX.F6wMTrv2'
$ws.Range("C14").Value = '$2b$10$B5qjGRTdetFzBfrGBaw4JevUd6hIWzmrD7tjk.iL2iOzA8e//yvkO'
$ws.Range("C15").Value = '$2b$10$.SPv/66kYKzuuyx7L.i84uRyoBSk.0VZ8Qbnn0u350Q4fxJGvyOw6'
$ws.Range("C16").Value = '$2b$10$W683IFxAM15CFlGUrY3g.uINpAQRnPzOWUBEzL9gGZKF9VO6p0A3e'
$ws.Range("C17").Value = '$2b$10$Mo22eEAO97GbfdNrXqG7POu1odE4DlkOnUzixtFpSAWkNOn8.izEm'
$ws.Range("C18").Value = '$2b$10$wO8cDxYezTX3wqWQlFUOvOYR/013Kb147UkLXqpw/mdgHRw41gZ02'
$ws.Range("C19").Value = '$2b$10$kyuwAHsKAWO.WFy5hFwtveY89/FrwgwDjO5flbQrAZk0Agdazfz/m'
$ws.Range("C20").Value = '$2b$10$8Xypu2JekmWPaqCr9CMicuCZ0MDYp6/h4gVGvBXaNyLdS9viwlpcu'
$ws.Range("C21").Value = '$2b$10$RYJZZxaVbMua1ZzqKjJqlOCVinyjMHVMRk7JHACkEzQswFTge1M02'
$ws.Range("C22").Value = '$2b$10$KRuMDcWfTCePVMszpo98cuX3npdk3N8ObnOn1OcKChm6FW8hRiuH2'
$ws.Range("C23").Value = '$2b$10$rhvcXl..PGYHZvOhgUplheQTzKXwjPYaebcbUHPGJ7mdspx0I1coa'
$ws.Range("C24").Value = '$2b$10$UIVD285dkqllY13KWt7ElOHqULYEoQf/XY5e.tmUHI9cbJYDOIeLi'
$ws.Range("C25").Value = '$2b$10$3pS.cpIG3Sk/s0r1g0s58.znQ3u/6kqQ.eERH5Zuq6uDNbXPvtdkK'
$ws.Range("C26").Value = '$2b$10$h8SbJcQ8NY2IuCag67Wyk.f7Ye173gbswSOunjpvtkqcDjmWm4qdG'
$ws.Range("C27").Value = '$2b$10$8w.WfWA0bmondadVzKJEkuNSpsy23.Cm6Jbt0nB0fARxyL.4zizfm'
$ws.Range("C28").Value = '$2b$10$JlWKzE4insGc.s082HLg2OxjUpzQ5s8m1tXa0Rhb.xcC8Kp5YZorK'
$ws.Range("C29").Value = '$2b$10$4bd9X4Z0YVIb4ltox1NVjuAoEViSj65R0jniASWwisaIWHL0XkyZu'
$ws.Range("C30").Value = '$2b$10$k0tXpN97pQ/iQEK5SwxUpuuKE9n2xa0ono4aXhEViZM3VA4sbMEdS'
$ws.Range("C31").Value = '$2b$10$9JnTrERVEyKcjQOhsA4rSe7Nd3x4P5R5xk2h8omvF9izxl6fBcBKC'
